$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "2020" column (X) to the table, mirroring the formatting of the
# existing "2019" column (W), then fill in the reported values.

$xlPasteFormats = -4122

$ws.Range("W4:W16").Copy()
$ws.Range("X4:X16").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("X4").Value2 = 2020
$ws.Range("X5").Value2 = 45.3
$ws.Range("X6").Value2 = 48.2
$ws.Range("X7").Value2 = 43.6
$ws.Range("X8").Value2 = 48.8
$ws.Range("X9").Value2 = 41.5
$ws.Range("X10").Value2 = 49.7
$ws.Range("X11").Value2 = 46.7
$ws.Range("X12").Value2 = 36.5
$ws.Range("X13").Value2 = 29.6
$ws.Range("X14").Value2 = 54.7
$ws.Range("X15").Value2 = 51.6
$ws.Range("X16").Value2 = 47.2

# Match the saved selection state recorded in the edited workbook.
[void]$ws.Range("AI21").Select()
